{"js": "const replacements = [\n  [\"939\u00f79=104, 3\", \"640\u00f77=91, 3\"],\n  [\"689\u00f75=137, 4\", \"297\u00f78=37, 1\"],\n  [\"561\u00f79=62, 3\", \"399\u00f79=44, 3\"],\n  [\"821\u00f72=410, 1\", \"290\u00f75=58, 0\"],\n  [\"556\u00f78=69, 4\", \"732\u00f79=81, 3\"],\n  [\"184\u00f76=30, 4\", \"542\u00f76=90, 2\"],\n  [\"965\u00f76=160, 5\", \"254\u00f73=84, 2\"],\n  [\"998\u00f76=166, 2\", \"912\u00f75=182, 2\"],\n  [\"269\u00f72=134, 1\", \"103\u00f78=12, 7\"],\n  [\"175\u00f74=43, 3\", \"132\u00f73=44, 0\"],\n  [\"377\u00f76=62, 5\", \"654\u00f79=72, 6\"],\n  [\"973\u00f74=243, 1\", \"755\u00f72=377, 1\"],\n  [\"170\u00f79=18, 8\", \"638\u00f72=319, 0\"],\n  [\"514\u00f78=64, 2\", \"594\u00f77=84, 6\"],\n  [\"476\u00f76=79, 2\", \"885\u00f77=126, 3\"],\n  [\"571\u00f74=142, 3\", \"171\u00f72=85, 1\"],\n  [\"130\u00f78=16, 2\", \"335\u00f73=111, 2\"],\n  [\"839\u00f79=93, 2\", \"814\u00f79=90, 4\"],\n  [\"718\u00f73=239, 1\", \"769\u00f72=384, 1\"],\n  [\"403\u00f77=57, 4\", \"512\u00f77=73, 1\"],\n  [\"824\u00f79=91, 5\", \"537\u00f73=179, 0\"],\n  [\"583\u00f72=291, 1\", \"216\u00f79=24, 0\"],\n  [\"714\u00f72=357, 0\", \"634\u00f77=90, 4\"],\n  [\"913\u00f73=304, 1\", \"700\u00f78=87, 4\"],\n  [\"215\u00f78=26, 7\", \"737\u00f77=105, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"939\u00f79=104, 3\", \"640\u00f77=91, 3\"),\n    @(\"689\u00f75=137, 4\", \"297\u00f78=37, 1\"),\n    @(\"561\u00f79=62, 3\", \"399\u00f79=44, 3\"),\n    @(\"821\u00f72=410, 1\", \"290\u00f75=58, 0\"),\n    @(\"556\u00f78=69, 4\", \"732\u00f79=81, 3\"),\n    @(\"184\u00f76=30, 4\", \"542\u00f76=90, 2\"),\n    @(\"965\u00f76=160, 5\", \"254\u00f73=84, 2\"),\n    @(\"998\u00f76=166, 2\", \"912\u00f75=182, 2\"),\n    @(\"269\u00f72=134, 1\", \"103\u00f78=12, 7\"),\n    @(\"175\u00f74=43, 3\", \"132\u00f73=44, 0\"),\n    @(\"377\u00f76=62, 5\", \"654\u00f79=72, 6\"),\n    @(\"973\u00f74=243, 1\", \"755\u00f72=377, 1\"),\n    @(\"170\u00f79=18, 8\", \"638\u00f72=319, 0\"),\n    @(\"514\u00f78=64, 2\", \"594\u00f77=84, 6\"),\n    @(\"476\u00f76=79, 2\", \"885\u00f77=126, 3\"),\n    @(\"571\u00f74=142, 3\", \"171\u00f72=85, 1\"),\n    @(\"130\u00f78=16, 2\", \"335\u00f73=111, 2\"),\n    @(\"839\u00f79=93, 2\", \"814\u00f79=90, 4\"),\n    @(\"718\u00f73=239, 1\", \"769\u00f72=384, 1\"),\n    @(\"403\u00f77=57, 4\", \"512\u00f77=73, 1\"),\n    @(\"824\u00f79=91, 5\", \"537\u00f73=179, 0\"),\n    @(\"583\u00f72=291, 1\", \"216\u00f79=24, 0\"),\n    @(\"714\u00f72=357, 0\", \"634\u00f77=90, 4\"),\n    @(\"913\u00f73=304, 1\", \"700\u00f78=87, 4\"),\n    @(\"215\u00f78=26, 7\", \"737\u00f77=105, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
